# Trade #59 (row 88 on "All Trades" / row 55 on "MarketMaking") closes at
# 2026-02-17 21:10:46 as an early_exit with a small gain, and a brand new
# trade #120 opens at 21:10:40 (rows appended to both sheets). Summary and
# Strategy Status sheets are refreshed to reflect the new aggregate stats.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Summary sheet
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.97
$summary.Range("B4").Value = 0.77
$summary.Range("B5").Value = 0.18
$summary.Range("B6").Value = 87
$summary.Range("B7").Value = 41
$summary.Range("B9").Value = 47.13

# ----------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row = row 5)
# ----------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.97
$status.Range("D5").Value = 54
$status.Range("E5").Value = 0.66
$status.Range("F5").Value = 0.97
$status.Range("G5").Value = 50

# ----------------------------------------------------------------------
# All Trades sheet - close out trade #87 (row 88)
# ----------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G88").Value = 0.79
$allTrades.Range("H88").Value = "CLOSED"
$allTrades.Range("I88").Value = 5.3333
$allTrades.Range("J88").Value = 0.04
$allTrades.Range("K88").Value = 100.97
$allTrades.Range("L88").Value = "early_exit"
$allTrades.Range("M88").Value = 0.14

# Append new trade #120 as row 121
$allTrades.Range("A121").Value = 120
$allTrades.Range("B121").NumberFormat = "@"
$allTrades.Range("B121").Value = "2026-02-17"
$allTrades.Range("C121").Value = "21:10:40"
$allTrades.Range("D121").Value = "MarketMaking"
$allTrades.Range("E121").Value = "UP"
$allTrades.Range("F121").Value = 0.75
$allTrades.Range("H121").Value = "OPEN"
$allTrades.Range("I121").Value = 0
$allTrades.Range("J121").Value = 0
$allTrades.Range("K121").Value = 100.9346450978375
$allTrades.Range("M121").Value = 0
$allTrades.Range("N121").Value = 0
$allTrades.Range("O121").Value = 0
$allTrades.Range("P121").Value = 0.6
$allTrades.Range("Q121").Value = "Normal spread capture: 19600 bps"

# ----------------------------------------------------------------------
# MarketMaking sheet - close out trade #87 (row 55)
# ----------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G55").Value = 0.79
$mm.Range("H55").Value = "CLOSED"
$mm.Range("I55").Value = 5.3333
$mm.Range("J55").Value = 0.04
$mm.Range("K55").Value = 100.97
$mm.Range("P55").Value = "early_exit"
$mm.Range("Q55").Value = 0.14

# Append new trade #120 as row 88
$mm.Range("A88").Value = 120
$mm.Range("B88").NumberFormat = "@"
$mm.Range("B88").Value = "2026-02-17"
$mm.Range("C88").Value = "21:10:40"
$mm.Range("D88").Value = "MarketMaking"
$mm.Range("E88").Value = "UP"
$mm.Range("F88").Value = 0.75
$mm.Range("H88").Value = "OPEN"
$mm.Range("I88").Value = 0
$mm.Range("J88").Value = 0
$mm.Range("K88").Value = 100.9346450978375
$mm.Range("L88").Value = 0
$mm.Range("M88").Value = 0
$mm.Range("N88").Value = 0.6
$mm.Range("O88").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q88").Value = 0
